$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Copy($ws.Range("A11"))
$ws.Range("A11").Value = "2021年"
$ws.Range("B11").Value = 8459.48
$ws.Range("C11").Value = 2020.93
$ws.Range("D11").Value = 417.91
$ws.Range("E11").Value = 0.93
$ws.Range("F11").Value = 5895.13
$ws.Range("G11").Value = 8177.76
$ws.Range("H11").Value = 1044.99
$ws.Range("I11").Value = 5111.17
$ws.Range("J11").Value = 782.01
$ws.Range("K11").Value = 993.23
$ws.Range("L11").Value = 767.61
$ws.Range("M11").Value = 80.34
$ws.Range("N11").Value = 2416.38
$ws.Range("O11").Value = 6930.65
$ws.Range("P11").Value = 323.28
$ws.Range("Q11").Value = 914.77
$ws.Range("R11").Value = 3568.44
$ws.Range("S11").Value = 254.59
$ws.Range("T11").Value = 8632.459999999999
$ws.Range("U11").Value = 4004.67
$ws.Range("V11").Value = 1286.72
$ws.Range("W11").Value = 423.39
$ws.Range("X11").Value = 2294.06
$ws.Range("Y11").Value = 11205.47
$ws.Range("Z11").Value = 1045.4
$ws.Range("AA11").Value = 5598.3
$ws.Range("AB11").Value = 184.3
$ws.Range("AC11").Value = 3550.48
$ws.Range("AD11").Value = 1880.49
$ws.Range("AE11").Value = 151021.62
$ws.Range("AF11").Value = 19387.19
$ws.Range("AG11").Value = 9048.120000000001
$ws.Range("AH11").Value = 1692.71
$ws.Range("AI11").Value = 3627.63
$ws.Range("AJ11").Value = 301.55
$ws.Range("AK11").Value = 5964.59
$ws.Range("AL11").Value = 5148.26
$ws.Range("AM11").Value = 6660.96
$ws.Range("AN11").Value = 320.26
$ws.Range("AO11").Value = 2177.9
$ws.Range("AP11").Value = 8045.95
$ws.Range("AQ11").Value = 381.21

$ws.Range("A10").Copy($ws.Range("A12"))
$ws.Range("A12").Value = "2022年"
$ws.Range("B12").Value = 9213.799999999999
$ws.Range("C12").Value = 2204.2
$ws.Range("D12").Value = 253.9
$ws.Range("E12").Value = 0.6
$ws.Range("F12").Value = 6546.5
$ws.Range("G12").Value = 8957.799999999999
$ws.Range("H12").Value = 1132
$ws.Range("I12").Value = 5414.8
$ws.Range("J12").Value = 801.1
$ws.Range("K12").Value = 951.2
$ws.Range("L12").Value = 1009.9
$ws.Range("M12").Value = 88.7
$ws.Range("N12").Value = 2469.8
$ws.Range("O12").Value = 7685
$ws.Range("P12").Value = 325.8
$ws.Range("Q12").Value = 1004
$ws.Range("R12").Value = 3741.6
$ws.Range("S12").Value = 277.7
$ws.Range("T12").Value = 9401.200000000001
$ws.Range("U12").Value = 4167.9
$ws.Range("V12").Value = 1442.1
$ws.Range("W12").Value = 418.8
$ws.Range("X12").Value = 2160.9
$ws.Range("Y12").Value = 12915.1
$ws.Range("Z12").Value = 1034.4
$ws.Range("AA12").Value = 6620.3
$ws.Range("AB12").Value = 195.4
$ws.Range("AC12").Value = 3600.1
$ws.Range("AD12").Value = 1875.9
$ws.Range("AE12").Value = 159134.4
$ws.Range("AF12").Value = 20736.8
$ws.Range("AG12").Value = 9283.200000000001
$ws.Range("AH12").Value = 1805.8
$ws.Range("AI12").Value = 4184.8
$ws.Range("AJ12").Value = 322.5
$ws.Range("AK12").Value = 5748.4
$ws.Range("AL12").Value = 2950.7
$ws.Range("AM12").Value = 7337.9
$ws.Range("AN12").Value = 376.4
$ws.Range("AO12").Value = 2471.6
$ws.Range("AP12").Value = 7602.4
$ws.Range("AQ12").Value = 403.5

